$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17.701128005981445
$ws.Range("C2").Value = 5.344827651977539
$ws.Range("D2").Value = 11.570488929748535
$ws.Range("E2").Value = 45.71428680419922
